$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 2298.75
$ws.Range("I41").Value = 468.2
$ws.Range("K41").Value = 468.2
$ws.Range("M41").Value = -28.19999999999999
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H76").Value = 1023.125
$ws.Range("I76").Value = 740.1429000000001
$ws.Range("K76").Value = 740.1429000000001
$ws.Range("M76").Value = -425.1429000000001
$ws.Range("H79").Value = 1023.125
$ws.Range("I79").Value = 740.1429000000001
$ws.Range("K79").Value = 740.1429000000001
$ws.Range("M79").Value = 351.8570999999999
$ws.Range("H80").Value = 466.4
$ws.Range("J80").Value = 498.33334
$ws.Range("L80").Value = 1495.00002
$ws.Range("N80").Value = -3491.00002
$ws.Range("H83").Value = 466.4
$ws.Range("J83").Value = 498.33334
$ws.Range("L83").Value = 4485.00006
$ws.Range("N83").Value = -14469.00006
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("H107").Value = 1362.7693
$ws.Range("I107").Value = 1459.75
$ws.Range("J107").Value = 199
$ws.Range("K107").Value = 1459.75
$ws.Range("L107").Value = 199
$ws.Range("M107").Value = 460.25
$ws.Range("N107").Value = -4039
$ws.Range("H113").Value = 4008
$ws.Range("I113").Value = 4008
$ws.Range("K113").Value = 4008
$ws.Range("M113").Value = -754
$ws.Range("H116").Value = 3934
$ws.Range("I116").Value = 3934
$ws.Range("K116").Value = 3934
$ws.Range("M116").Value = -492
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 9808.73
$ws.Range("I132").Value = 13001.883
$ws.Range("K132").Value = 39005.649
$ws.Range("M132").Value = -36475.649
$ws.Range("H138").Value = 3086.8333
$ws.Range("I138").Value = 1434.6666
$ws.Range("K138").Value = 4303.9998
$ws.Range("M138").Value = 836.0002000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H45").Value = 1898.381
$ws.Range("I45").Value = 1108
$ws.Range("K45").Value = 1108
$ws.Range("M45").Value = -731
$ws.Range("H61").Value = 7399.8
$ws.Range("I61").Value = 6333
$ws.Range("K61").Value = 6333
$ws.Range("M61").Value = -6121
$ws.Range("H130").Value = 18143
$ws.Range("J130").Value = 18143
$ws.Range("L130").Value = 18143
$ws.Range("N130").Value = -28183
$ws.Range("H136").Value = 7399.8
$ws.Range("I136").Value = 6333
$ws.Range("K136").Value = 18999
$ws.Range("M136").Value = -16449

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2808.3125
$ws.Range("I20").Value = 1664.1666
$ws.Range("K20").Value = 1664.1666
$ws.Range("M20").Value = -1417.1666
$ws.Range("H42").Value = 240000
$ws.Range("J42").Value = 240000
$ws.Range("L42").Value = 240000
$ws.Range("N42").Value = -240656
$ws.Range("H107").Value = 4217.737
$ws.Range("I107").Value = 2768.6
$ws.Range("J107").Value = 5827.8887
$ws.Range("K107").Value = 2768.6
$ws.Range("L107").Value = 5827.8887
$ws.Range("M107").Value = -848.5999999999999
$ws.Range("N107").Value = -9667.8887

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5259.475
$ws.Range("I31").Value = 3110.7778
$ws.Range("J31").Value = 9722.154
$ws.Range("K31").Value = 3110.7778
$ws.Range("L31").Value = 9722.154
$ws.Range("M31").Value = -2815.7778
$ws.Range("N31").Value = -10312.154
$ws.Range("H34").Value = 5259.475
$ws.Range("I34").Value = 3110.7778
$ws.Range("J34").Value = 9722.154
$ws.Range("K34").Value = 3110.7778
$ws.Range("L34").Value = 9722.154
$ws.Range("M34").Value = -2908.7778
$ws.Range("N34").Value = -10126.154
$ws.Range("H122").Value = 809.5333000000001
$ws.Range("I122").Value = 809.5333000000001
$ws.Range("K122").Value = 2428.5999
$ws.Range("M122").Value = 21.40009999999984
$ws.Range("H134").Value = 1727.5
$ws.Range("I134").Value = 1810.2
$ws.Range("K134").Value = 5430.6
$ws.Range("M134").Value = -2895.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5000
$ws.Range("I3").Value = 5000
$ws.Range("K3").Value = 15000
$ws.Range("M3").Value = -14888
$ws.Range("H109").Value = 999.5
$ws.Range("I109").Value = 999.5
$ws.Range("K109").Value = 2998.5
$ws.Range("M109").Value = -1958.5
$ws.Range("H110").Value = 19999
$ws.Range("I110").Value = 19999
$ws.Range("K110").Value = 59997
$ws.Range("M110").Value = -55907
$ws.Range("H113").Value = 1279.1875
$ws.Range("I113").Value = 1000.8333
$ws.Range("J113").Value = 1446.2
$ws.Range("K113").Value = 3002.4999
$ws.Range("L113").Value = 4338.6
$ws.Range("M113").Value = -832.4998999999998
$ws.Range("N113").Value = -8678.6
$ws.Range("H121").Value = 1158.4286
$ws.Range("I121").Value = 527.25
$ws.Range("K121").Value = 1581.75
$ws.Range("M121").Value = -271.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 28000
$ws.Range("J14").Value = 1000
$ws.Range("L14").Value = 1000
$ws.Range("N14").Value = -1336
$ws.Range("H32").Value = 64000
$ws.Range("J32").Value = 64000
$ws.Range("L32").Value = 64000
$ws.Range("N32").Value = -64592
$ws.Range("H102").Value = 2094.5334
$ws.Range("J102").Value = 4102
$ws.Range("L102").Value = 4102
$ws.Range("N102").Value = -7346
$ws.Range("H113").Value = 3751.95
$ws.Range("I113").Value = 2629
$ws.Range("J113").Value = 7120.8
$ws.Range("K113").Value = 2629
$ws.Range("L113").Value = 7120.8
$ws.Range("M113").Value = -459
$ws.Range("N113").Value = -11460.8
$ws.Range("H122").Value = 3852.7693
$ws.Range("I122").Value = 2787.3333
$ws.Range("K122").Value = 8361.999899999999
$ws.Range("M122").Value = -5911.999899999999
$ws.Range("H126").Value = 1687.5
$ws.Range("I126").Value = 1500
$ws.Range("J126").Value = 1875
$ws.Range("K126").Value = 4500
$ws.Range("L126").Value = 5625
$ws.Range("M126").Value = -2030
$ws.Range("N126").Value = -10565
$ws.Range("H141").Value = 70000
$ws.Range("J141").Value = 70000
$ws.Range("L141").Value = 70000
$ws.Range("N141").Value = -80360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 214.75
$ws.Range("I9").Value = 81.333336
$ws.Range("K9").Value = 81.333336
$ws.Range("M9").Value = 142.666664
$ws.Range("H10").Value = 2366
$ws.Range("I10").Value = 2799.5
$ws.Range("J10").Value = 1499
$ws.Range("K10").Value = 2799.5
$ws.Range("L10").Value = 1499
$ws.Range("M10").Value = -2659.5
$ws.Range("N10").Value = -1779
$ws.Range("H12").Value = 651.5
$ws.Range("I12").Value = 203
$ws.Range("J12").Value = 801
$ws.Range("K12").Value = 203
$ws.Range("L12").Value = 801
$ws.Range("M12").Value = -33
$ws.Range("N12").Value = -1141
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("I19").Value = 500
$ws.Range("J19").Value = 5500
$ws.Range("K19").Value = 500
$ws.Range("L19").Value = 5500
$ws.Range("M19").Value = -330
$ws.Range("N19").Value = -5840
$ws.Range("H46").Value = 5937.375
$ws.Range("I46").Value = 499.66666
$ws.Range("K46").Value = 499.66666
$ws.Range("M46").Value = -311.66666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 59999.668
$ws.Range("I27").Value = 59999
$ws.Range("K27").Value = 59999
$ws.Range("M27").Value = -59930
$ws.Range("H62").Value = 9556.111000000001
$ws.Range("I62").Value = 7000.5
$ws.Range("K62").Value = 7000.5
$ws.Range("M62").Value = -6376.5
$ws.Range("H65").Value = 9556.111000000001
$ws.Range("I65").Value = 7000.5
$ws.Range("K65").Value = 35002.5
$ws.Range("M65").Value = -31882.5
$ws.Range("H124").Value = 100000
$ws.Range("J124").Value = 100000
$ws.Range("L124").Value = 100000
$ws.Range("N124").Value = -109820
